# Update results with new shocks:
#  - Row 1 (C1:Q1) on both "Corn" and "Soybean" sheets used to hold the
#    shared-string labels r_yr1..r_yr15; they become plain numbered
#    columns (1..15) built from live formulas (C1=1, D1=C1+1, and a
#    shared formula E1:Q1 = previous cell + 1).
#  - Refresh the sheet/window selection state left behind by the edit.

$wb = $excel.ActiveWorkbook

$wsCorn = $wb.Worksheets.Item("Corn")
$wsSoy  = $wb.Worksheets.Item("Soybean")

foreach ($ws in @($wsCorn, $wsSoy)) {
    $ws.Range("C1").Value = 1
    $ws.Range("D1").Formula = "=C1+1"
    $ws.Range("E1:Q1").Formula = "=D1+1"
}

# Soybean: no longer the selected tab; scrolled so column B is the
# left-most visible column, with C1:Q1 selected.
$wsSoy.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$wsSoy.Range("C1:Q1").Select() | Out-Null

# Corn: becomes the selected tab, with C13 as the active cell.
$wsCorn.Activate() | Out-Null
$wsCorn.Range("C13").Select() | Out-Null
